$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.124.29"
$ws.Range("E2").Value = "  -3.06%  "
$ws.Range("D3").Value = "1.867.63"
$ws.Range("E3").Value = "  -2.40%  "
$ws.Range("E4").Value = "  +0.15%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "307.13"
$c.Style = "Normal"
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +0.22%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.5121"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +2.10%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3736"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -2.07%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.07137"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -2.48%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.8876"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -2.77%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "20.64"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -3.01%  "
$ws.Range("D12").Value = "1.869.05"
$ws.Range("E12").Value = "  -2.49%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.07537"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -1.67%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "5.313"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -3.38%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "89.12"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -4.01%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +0.13%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.000008471"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -3.16%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "14.11"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -3.93%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +0.16%  "
$ws.Range("D20").Value = "27.203.77"
$ws.Range("E20").Value = "  -2.94%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "5.047"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -2.81%  "
$ws.Range("D22").Value = "2.093.76"
$ws.Range("E22").Value = "  -1.53%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "10.55"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -2.91%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "6.472"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -2.37%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "150.15"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -2.12%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "1.838"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -0.01%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "17.91"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -2.96%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "2.091"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -5.84%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "112.67"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -2.51%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "4.745"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -3.84%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "4.669"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -3.89%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.09009"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -0.45%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.05127"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -2.92%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "3.100"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -3.43%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "1.159"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -6.40%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.7333"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -6.16%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.02045"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -2.05%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "2.499"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -5.31%  "
$ws.Range("E39").Value = "  -1.01%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "1.075"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -1.81%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.5305"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -4.63%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "6.604"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -3.85%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "116.52"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +1.97%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "8.337"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -2.64%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.1472"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -3.38%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "1.000"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +0.20%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.4616"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -4.54%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "9.966"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -6.28%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "1.568"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -4.41%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "64.53"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -4.64%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "36.45"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -1.88%  "
